# Added negative test cases for the GKuser
# Populate rows 18-22 (columns A/B) on the active sheet with new
# negative-test-case key/value pairs. Cells are written in row-major
# order (A18,B18,A19,B19,...) so that shared strings are interned in
# the same order the authoring tool produced them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = 'invalid.username'
$ws.Range("B18").Value = '63654,$%^&**,hdsj,  ,'

$ws.Range("A19").Value = 'invalid.answer'
$ws.Range("B19").Value = '63654,$%^&**,hdsj,  ,'

$ws.Range("A20").Value = 'unique.Username'
$ws.Range("B20").Value = '1,12,123,123456789qwertsdjksjdhjfsdjkdsfjfjkfjskdsjjjfdsjkfsdjkdsjkfhjfksdhjfksdhfjksd,   ,'

$ws.Range("A21").Value = 'invalid.username'
$ws.Range("B21").Value = '63654,$%^&**,hdsj,  ,'

$ws.Range("A22").Value = 'invalid.password'
$ws.Range("B22").Value = 'password,123456,myname123,1234567890,aaaaaaaa'

# Move the selection to the last cell entered, matching where the
# author's cursor ended up after typing in the new rows.
$ws.Range("B22").Select() | Out-Null
